# Updated cryptos list - refresh Price and Volume(1h) figures scraped from coinranking.com
# (also corrects the row order for Elrond/EnergySwap which had swapped positions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value2 = "30.164.09"
$ws.Range("E2").Value2 = "  -0.91%  "
# Row 3: Ethereum
$ws.Range("D3").Value2 = "1.847.60"
$ws.Range("E3").Value2 = "  -2.32%  "
# Row 4: TetherUSD
$ws.Range("E4").Value2 = "  -0.12%  "
# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "235.97"
$ws.Range("E5").Value2 = "  -0.92%  "
# Row 6: USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "1.0000"
$ws.Range("E6").Value2 = "  -0.13%  "
# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4781"
$ws.Range("E7").Value2 = "  -2.48%  "
# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.2802"
$ws.Range("E8").Value2 = "  -4.55%  "
# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06470"
$ws.Range("E9").Value2 = "  -3.43%  "
# Row 10: WrappedEther
$ws.Range("D10").Value2 = "1.856.57"
$ws.Range("E10").Value2 = "  -1.70%  "
# Row 11: TRON
$ws.Range("E11").Value2 = "  -0.49%  "
# Row 12: Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "16.24"
$ws.Range("E12").Value2 = "  -4.27%  "
# Row 13: Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "5.101"
$ws.Range("E13").Value2 = "  -0.52%  "
# Row 14: Litecoin
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "87.03"
$ws.Range("E14").Value2 = "  -0.83%  "
# Row 15: Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.6444"
$ws.Range("E15").Value2 = "  -2.87%  "
# Row 16: WrappedBTC
$ws.Range("D16").Value2 = "30.094.04"
$ws.Range("E16").Value2 = "  -1.07%  "
# Row 17: Dai
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "1.000"
$ws.Range("E17").Value2 = "  -0.01%  "
# Row 18: Avalanche
$ws.Range("E18").Value2 = "  -1.78%  "
# Row 19: ShibaInu
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.000007618"
$ws.Range("E19").Value2 = "  -2.74%  "
# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "224.75"
$ws.Range("E20").Value2 = "  +18.20%  "
# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value2 = "2.098.39"
$ws.Range("E21").Value2 = "  -1.57%  "
# Row 22: BinanceUSD
$ws.Range("E22").Value2 = "  -0.02%  "
# Row 23: Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.278"
# Row 24: Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "6.066"
$ws.Range("E24").Value2 = "  -1.02%  "
# Row 25: Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "9.196"
$ws.Range("E25").Value2 = "  -2.90%  "
# Row 26: Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "163.47"
$ws.Range("E26").Value2 = "  -0.13%  "
# Row 27: EthereumClassic
$ws.Range("E27").Value2 = "  +1.34%  "
# Row 28: LidoDAOToken
$ws.Range("E28").Value2 = "  -0.71%  "
# Row 29: Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.427"
$ws.Range("E29").Value2 = "  -2.73%  "
# Row 30: Stellar
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "0.09189"
$ws.Range("E30").Value2 = "  +0.45%  "
# Row 31: InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "4.232"
$ws.Range("E31").Value2 = "  -2.95%  "
# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "3.948"
$ws.Range("E32").Value2 = "  -2.24%  "
# Row 33: Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.05008"
$ws.Range("E33").Value2 = "  -3.70%  "
# Row 34: ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.7374"
$ws.Range("E34").Value2 = "  -0.33%  "
# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.134"
# Row 36: HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "2.687"
$ws.Range("E36").Value2 = "  -1.08%  "
# Row 37: VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.01810"
$ws.Range("E37").Value2 = "  -0.26%  "
# Row 38: MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.601"
$ws.Range("E38").Value2 = "  -2.74%  "
# Row 39: TrustWalletToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.9041"
$ws.Range("E39").Value2 = "  -1.83%  "
# Row 40: RenderToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "2.053"
$ws.Range("E40").Value2 = "  +0.79%  "
# Row 41: FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "5.942"
$ws.Range("E41").Value2 = "  +0.07%  "
# Row 42: Quant
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "106.80"
$ws.Range("E42").Value2 = "  +0.72%  "
# Row 43: TheSandbox
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.4242"
$ws.Range("E43").Value2 = "  -3.56%  "
# Row 44: PaxDollar
$ws.Range("E44").Value2 = "  +0.59%  "
# Row 45: Aptos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "7.376"
$ws.Range("E45").Value2 = "  -2.56%  "
# Row 46: Algorand
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.1316"
$ws.Range("E46").Value2 = "  -4.00%  "
# Row 47: NEARProtocol
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "1.557"
$ws.Range("E47").Value2 = "  +10.85%  "
# Row 48: Aave
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "63.93"
$ws.Range("E48").Value2 = "  -6.26%  "
# Row 49: EnergySwap
$ws.Range("B49").Value2 = "Elrond"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "34.22"
$ws.Range("E49").Value2 = "  -2.14%  "
# Row 50: Elrond
$ws.Range("B50").Value2 = "EnergySwap"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "8.736"
$ws.Range("E50").Value2 = "  -2.89%  "
# Row 51: Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.05660"
$ws.Range("E51").Value2 = "  -2.76%  "
